$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.249.84"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "'1.861.15"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'235.56"

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "'0.4668"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").Value = "'0.2832"
$ws.Range("E8").Value = "  +0.68%  "

$ws.Range("D9").Value = "'0.06507"
$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").Value = "'21.37"
$ws.Range("E10").Value = "  +6.86%  "

$ws.Range("D11").Value = "'0.07910"
$ws.Range("E11").Value = "  +1.50%  "

$ws.Range("D12").Value = "'97.09"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "'1.865.31"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").Value = "'5.144"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "'0.6776"
$ws.Range("E15").Value = "  +1.91%  "

$ws.Range("D16").Value = "'278.22"
$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("D17").Value = "'30.251.52"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "'13.70"
$ws.Range("E18").Value = "  +8.95%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").Value = "'5.382"
$ws.Range("E20").Value = "  -1.42%  "

$ws.Range("D21").Value = "'2.109.87"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'0.000007307"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'167.13"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").Value = "'9.154"

$ws.Range("D27").Value = "'19.04"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").Value = "'1.920"
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("E29").Value = "  +2.82%  "

$ws.Range("D30").Value = "'0.09694"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").Value = "'4.362"
$ws.Range("E31").Value = "  -1.29%  "

$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("E33").Value = "  -1.64%  "

$ws.Range("D34").Value = "'0.04719"
$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").Value = "'1.126"
$ws.Range("E35").Value = "  +2.50%  "

$ws.Range("D36").Value = "'0.7036"
$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("D37").Value = "'2.710"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'0.01859"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").Value = "'2.590"
$ws.Range("E39").Value = "  +3.41%  "

$ws.Range("D40").Value = "'6.330"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").Value = "'75.27"
$ws.Range("E41").Value = "  +4.59%  "

$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").Value = "'0.8477"
$ws.Range("E43").Value = "  -1.24%  "

$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("D46").Value = "'103.36"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").Value = "'972.00"
$ws.Range("E47").Value = "  -4.20%  "

$ws.Range("D48").Value = "'7.142"
$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("D49").Value = "'9.263"
$ws.Range("E49").Value = "  +2.92%  "

$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1127"
$ws.Range("E51").Value = "  -1.19%  "
